$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Wk[12] Sunday 3.6.18"
$ws.Range("B13").Value = "1200 - 1200"
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = "Debugging blocking logic."

$ws.Range("A14").Select()
